$d = $word.ActiveDocument

# --- Step 1: fix the visible text for the "experienced with" -> "with experience in" rewrite ---
$found = $d.Content.Find.Execute(
    "experienced with", $true, $false, $false, $false, $false,
    $true, 1, $false, "with experience in", 2)

# --- Step 2: low-level XML surgery for the run-split / proofErr markup that Word's COM
#     text-editing primitives can't express (run boundaries + spell-check proofErr tags) ---

$full = $d.Content
$pkg = $full.WordOpenXML

$partMarker = '<pkg:part pkg:name="/word/document.xml"'
$partIdx = $pkg.IndexOf($partMarker)
$dataOpen = '<pkg:xmlData>'
$dataClose = '</pkg:xmlData>'
$dataStart = $pkg.IndexOf($dataOpen, $partIdx) + $dataOpen.Length
$dataEnd = $pkg.IndexOf($dataClose, $dataStart)
$docXml = $pkg.Substring($dataStart, $dataEnd - $dataStart)

# --- 2a: "... graduate, with experience in object-oriented programming and full-stack
#          development. Pursuing full-time opportunities in the industry." currently
#          lives in one merged run (same rPr throughout) -- split it back into the
#          runs the target markup expects (only the "graduate, ... development" part
#          actually changed; the rest is restored to its original run boundaries). ---
$oldChunk = '<w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>graduate, with experience in object-oriented programming and full-stack development. Pursuing full-time opportunities in the industry.</w:t></w:r>'
$rPr = '<w:rPr><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr>'
$newChunk = '<w:r w:rsidR="001F11E7">' + $rPr + '<w:t>graduate</w:t></w:r>' + `
    '<w:r w:rsidRPr="00112E3E">' + $rPr + '<w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">with experience in </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>object-oriented programming and full-stack development</w:t></w:r>' + `
    '<w:r w:rsidRPr="00112E3E">' + $rPr + '<w:t>. Pursuing</w:t></w:r>' + `
    '<w:r w:rsidR="000C1F1B">' + $rPr + '<w:t xml:space="preserve"> full-time</w:t></w:r>' + `
    '<w:r w:rsidRPr="00112E3E">' + $rPr + '<w:t xml:space="preserve"> opportunities in the industry.</w:t></w:r>'
if ($docXml.IndexOf($oldChunk) -lt 0) {
    throw "chunk 2a not found"
}
$docXml = $docXml.Replace($oldChunk, $newChunk)

# --- 2b: ", Cybersec" -> ", " + proofErr-wrapped "Cybersec" run ---
$oldChunk = '<w:r w:rsidR="008325DC" w:rsidRPr="00841171"><w:rPr><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>, Cybersec</w:t></w:r>'
$rPr2 = '<w:rPr><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr>'
$newChunk = '<w:r w:rsidR="008325DC" w:rsidRPr="00841171">' + $rPr2 + '<w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r w:rsidR="008325DC" w:rsidRPr="00841171">' + $rPr2 + '<w:t>Cybersec</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
if ($docXml.IndexOf($oldChunk) -lt 0) {
    throw "chunk 2b not found"
}
$docXml = $docXml.Replace($oldChunk, $newChunk)

# --- 2c: wrap the bold "Firestore" run with proofErr spellStart/spellEnd ---
$oldChunk = '<w:r w:rsidRPr="00C27B7F"><w:rPr><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Firestore</w:t></w:r>'
$newChunk = '<w:proofErr w:type="spellStart"/>' + $oldChunk + '<w:proofErr w:type="spellEnd"/>'
if ($docXml.IndexOf($oldChunk) -lt 0) {
    throw "chunk 2c not found"
}
$docXml = $docXml.Replace($oldChunk, $newChunk)

$newPkg = $pkg.Substring(0, $dataStart) + $docXml + $pkg.Substring($dataEnd)
$full.WordOpenXML = $newPkg
